# EX_Reactions.xlsx - add the missing "EX_Fat(e)" transport reaction row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reaction was found - insert it as row 2 (just below the header),
# pushing every existing reaction down by one row.
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "EX_Fat_LPAREN_e_RPAREN_"

# Give the newly inserted row the same "just edited" look Excel applies
# (slightly shaded text, vertically centered) so it stands out from the rest.
$newCell = $ws.Range("A2")
$newCell.Font.ThemeColor = 1
$newCell.VerticalAlignment = -4108

# Reset the view: scroll back to the top and select F7 instead of the old
# F108 selection that pointed past the end of the (now longer) list.
$ws.Range("F7").Select() | Out-Null

# Make sure the sheet has a concrete page setup (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "Inserted EX_Fat_LPAREN_e_RPAREN_ at A2"
